$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B3").Value = 2
$ws.Range("E3").Value = 3
$ws.Range("H3").Value = 4
$ws.Range("K3").Value = 5
$ws.Range("B4").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("H4").Value = 2
$ws.Range("K4").Value = 5
$ws.Range("B5").Value = 3
$ws.Range("E5").Value = 4
$ws.Range("H5").Value = 2
$ws.Range("K5").Value = 5
$ws.Range("B6").Value = 3
$ws.Range("E6").Value = 5
$ws.Range("H6").Value = 2
$ws.Range("K6").Value = 5
$ws.Range("B7").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("H7").Value = 3
$ws.Range("K7").Value = 5
$ws.Range("B8").Value = 3
$ws.Range("E8").Value = 5
$ws.Range("H8").Value = 3
$ws.Range("K8").Value = 2
$ws.Range("B9").Value = 2
$ws.Range("E9").Value = 5
$ws.Range("H9").Value = 3
$ws.Range("K9").Value = 5
$ws.Range("B10").Value = 5
$ws.Range("E10").Value = 3
$ws.Range("H10").Value = 2
$ws.Range("K10").Value = 3
$ws.Range("B11").Value = 3
$ws.Range("E11").Value = 1
$ws.Range("H11").Value = 2
$ws.Range("K11").Value = 5
$ws.Range("B12").Value = 4
$ws.Range("E12").Value = 3
$ws.Range("H12").Value = 3
$ws.Range("K12").Value = 1
$ws.Range("B13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 2
$ws.Range("K13").Value = 5
$ws.Range("B14").Value = 4
$ws.Range("E14").Value = 3
$ws.Range("H14").Value = 3
$ws.Range("K14").Value = 2
$ws.Range("B15").Value = 3
$ws.Range("E15").Value = 3
$ws.Range("H15").Value = 1
$ws.Range("K15").Value = 2
$ws.Range("B16").Value = 3
$ws.Range("E16").Value = 2
$ws.Range("H16").Value = 2
$ws.Range("K16").Value = 5
$ws.Range("B17").Value = 4
$ws.Range("E17").Value = 4
$ws.Range("H17").Value = 3
$ws.Range("K17").Value = 1
$ws.Range("B18").Value = 4
$ws.Range("E18").Value = 4
$ws.Range("H18").Value = 4
$ws.Range("K18").Value = 2
$ws.Range("B19").Value = 3
$ws.Range("E19").Value = 1
$ws.Range("H19").Value = 5
$ws.Range("K19").Value = 4
$ws.Range("B20").Value = 3
$ws.Range("E20").Value = 4
$ws.Range("H20").Value = 4
$ws.Range("K20").Value = 5
$ws.Range("B21").Value = 2
$ws.Range("E21").Value = 2
$ws.Range("H21").Value = 4
$ws.Range("K21").Value = 5
$ws.Range("B22").Value = 1
$ws.Range("E22").Value = 4
$ws.Range("H22").Value = 1
$ws.Range("K22").Value = 4
$ws.Range("O19").Select()
